# "Updated labels on figures" - the header label in B1 was renamed
# from the Dutch "naam" to the English "name".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B1").Value = "name"
